$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1) text while keeping existing styles
$ws.Range("A1").Value = "Dolar Compra"
$ws.Range("B1").Value = "Dolar Venda"
$ws.Range("C1").Value = "Euro Compra"
$ws.Range("D1").Value = "Euro Venda"

# Copy style from an existing header cell (A1 already has style s=1) to new D1 cell
$ws.Range("D1").Style = $ws.Range("A1").Style

# Clear old row 3 (Euro row) entirely since new data only occupies rows 1-2
$ws.Range("A3:C3").Clear()

# Set row 2 values as text (inline strings) per new API-driven quote format
$ws.Range("A2").Value = "5.7991"
$ws.Range("B2").Value = "5.7991"
$ws.Range("C2").Value = "5.7991"
$ws.Range("D2").Value = "5.7991"
